$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Saturday")
$ws.Range("B16").Value = '(0, ''theory of computation'', ''small_lec'', ''5 csen ii'', 53)
'
$ws.Range("B17").Value = '(0, ''graphics'', ''lab'', ''5 dmet 23'', 63)
(0, ''data bases'', ''lab'', ''5 dmet 22'', 62)
'

$ws.Range("C16").Value = '(1, ''medi'', ''tut'', ''5 csen ii 18'', 38)
(1, ''dsd'', ''tut'', ''5 csen ii 17'', 37)
(1, ''comm ne'', ''tut'', ''5 csen ii 19'', 36)
(1, ''math'', ''tut'', ''5 csen ii 20'', 35)
(1, ''data bases'', ''lab'', ''5 csen ii 21'', 63)
'
$ws.Range("C17").Value = '(1, ''math'', ''tut'', ''5 dmet 23'', 34)
(1, ''data bases'', ''tut'', ''5 dmet 22'', 33)
'

$ws.Range("D16").Value = '(2, ''theory of computation'', ''tut'', ''5 csen ii 21'', 31)
(2, ''medi'', ''tut'', ''5 csen ii 19'', 30)
(2, ''math'', ''tut'', ''5 csen ii 16'', 29)
'
$ws.Range("D17").Value = '(2, ''computer graphics'', ''small_lec'', ''5 dmet'', 52)
'

$ws.Range("E16").Value = '(3, ''theory of computation'', ''tut'', ''5 csen ii 17'', 38)
(3, ''theory of computation'', ''tut'', ''5 csen ii 18'', 37)
(3, ''dsd'', ''tut'', ''5 csen ii 19'', 36)
(3, ''comm ne'', ''tut'', ''5 csen ii 21'', 35)
(3, ''medi'', ''tut'', ''5 csen ii 16'', 34)
(3, ''dsd'', ''tut'', ''5 csen ii 20'', 33)
'
$ws.Range("E17").Value = '(3, ''graphics'', ''tut'', ''5 dmet 23'', 32)
'

$ws = $wb.Worksheets.Item("Sunday")
$ws.Range("B16").Value = '(5, ''medi'', ''tut'', ''5 csen ii 17'', 36)
(5, ''data bases'', ''lab'', ''5 csen ii 20'', 61)
(5, ''data bases'', ''tut'', ''5 csen ii 18'', 35)
(5, ''data bases'', ''lab'', ''5 csen ii 19'', 60)
(5, ''dsd'', ''tut'', ''5 csen ii 21'', 34)
(5, ''theory of computation'', ''tut'', ''5 csen ii 16'', 33)
'
$ws.Range("B17").Value = '(5, ''comm ne'', ''tut'', ''5 dmet 23'', 32)
(5, ''graphics'', ''lab'', ''5 dmet 22'', 59)
'

$ws.Range("E16").Value = '(8, ''data bases'', ''lab'', ''5 csen ii 17'', 57)
(8, ''theory of computation'', ''tut'', ''5 csen ii 20'', 35)
(8, ''comm ne'', ''tut'', ''5 csen ii 18'', 34)
(8, ''theory of computation'', ''tut'', ''5 csen ii 19'', 33)
(8, ''comm ne'', ''tut'', ''5 csen ii 16'', 32)
'
$ws.Range("E17").Value = '(8, ''medi'', ''tut'', ''5 dmet 23'', 31)
(8, ''dsd'', ''tut'', ''5 dmet 22'', 30)
'

$ws = $wb.Worksheets.Item("Monday")
$ws.Range("B16").Value = '(10, ''theory of computation'', ''small_lec'', ''5 csen ii'', 52)
'
$ws.Range("B17").Value = '(10, ''data bases'', ''tut'', ''5 dmet 23'', 36)
(10, ''math'', ''tut'', ''5 dmet 22'', 35)
'

$ws.Range("D16").Value = '(12, ''comm ne'', ''tut'', ''5 csen ii 20'', 33)
(12, ''data bases'', ''tut'', ''5 csen ii 16'', 32)
(12, ''data bases'', ''tut'', ''5 csen ii 21'', 31)
(12, ''math'', ''tut'', ''5 csen ii 17'', 30)
(12, ''math'', ''tut'', ''5 csen ii 19'', 29)
(12, ''data bases'', ''lab'', ''5 csen ii 18'', 59)
'
$ws.Range("D17").Value = '(12, ''dsd'', ''tut'', ''5 dmet 23'', 28)
(12, ''medi'', ''tut'', ''5 dmet 22'', 27)
'

$ws.Range("E16").Value = '(13, ''data bases'', ''tut'', ''5 csen ii 20'', 27)
(13, ''dsd'', ''tut'', ''5 csen ii 16'', 26)
(13, ''medi'', ''tut'', ''5 csen ii 21'', 25)
(13, ''math'', ''tut'', ''5 csen ii 18'', 24)
(13, ''data bases'', ''tut'', ''5 csen ii 17'', 23)
'
$ws.Range("E17").Value = '(13, ''data bases'', ''lab'', ''5 dmet 23'', 56)
(13, ''comm ne'', ''tut'', ''5 dmet 22'', 22)
'
